$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.824.25"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.734.56"
$ws.Range("E3").Value = "  -0.67%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.16"
$ws.Range("E5").Value = "  -1.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5142"
$ws.Range("E7").Value = "  +1.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2790"
$ws.Range("E8").Value = "  +4.45%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.32"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06103"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").Value = "1.747.18"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07033"
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6433"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.520"
$ws.Range("E15").Value = "  +1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.74"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "25.825.44"
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006601"
$ws.Range("E21").Value = "  -0.52%  "
$ws.Range("D22").Value = "1.973.16"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.136"
$ws.Range("E23").Value = "  +2.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.667"
$ws.Range("E24").Value = "  +4.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.124"
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.05"
$ws.Range("E26").Value = "  +2.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.512"
$ws.Range("E27").Value = "  +3.25%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.06"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.797"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.10"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08324"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.678"
$ws.Range("E32").Value = "  +0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.420"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04487"
$ws.Range("E34").Value = "  +2.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.611"
$ws.Range("E35").Value = "  -1.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9816"
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6108"
$ws.Range("E37").Value = "  +1.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.648"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.936"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9998"
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "100.31"
$ws.Range("E42").Value = "  -1.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3829"
$ws.Range("E43").Value = "  +0.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7276"
$ws.Range("E44").Value = "  -2.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.963"
$ws.Range("E45").Value = "  +1.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05395"
$ws.Range("E46").Value = "  -1.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.280"
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1118"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "52.91"
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "29.98"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.600"
$ws.Range("E51").Value = "  +2.82%  "
